$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2291666666666667
$ws.Range("J2").Value = 0.0625
$ws.Range("S2").Value = 0.08333333333333333
# Row 3
$ws.Range("J3").Value = 0.1
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.1
# Row 4
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.4166666666666667
# Row 6
$ws.Range("B6").Value = 0.09677419354838709
$ws.Range("D6").Value = 0.06451612903225806
$ws.Range("F6").Value = 0.03225806451612903
$ws.Range("J6").Value = 0.2258064516129032
$ws.Range("Q6").Value = 0.1612903225806452
$ws.Range("R6").Value = 0.1290322580645161
$ws.Range("S6").Value = 0.2903225806451613
# Row 7
$ws.Range("D7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.2941176470588235
$ws.Range("Q7").Value = 0.05882352941176471
$ws.Range("S7").Value = 0.5882352941176471
# Row 8
$ws.Range("B8").Value = 0.06557377049180328
$ws.Range("D8").Value = 0.03278688524590164
$ws.Range("F8").Value = 0.04918032786885246
$ws.Range("J8").Value = 0.1967213114754098
$ws.Range("Q8").Value = 0.1639344262295082
$ws.Range("R8").Value = 0.09836065573770492
$ws.Range("S8").Value = 0.3934426229508197
# Row 9
$ws.Range("B9").Value = 0.03703703703703703
$ws.Range("D9").Value = 0.03703703703703703
$ws.Range("J9").Value = 0.1851851851851852
$ws.Range("O9").Value = 0.03703703703703703
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("R9").Value = 0.1111111111111111
$ws.Range("S9").Value = 0.3703703703703703
# Row 10
$ws.Range("B10").Value = 0.09246575342465753
$ws.Range("D10").Value = 0.02054794520547945
$ws.Range("F10").Value = 0.0684931506849315
$ws.Range("J10").Value = 0.1575342465753425
$ws.Range("O10").Value = 0.02054794520547945
$ws.Range("Q10").Value = 0.2636986301369863
$ws.Range("R10").Value = 0.1061643835616438
$ws.Range("S10").Value = 0.2705479452054795
# Row 11
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.2666666666666667
$ws.Range("L11").Value = 0.4666666666666667
# Row 12
$ws.Range("J12").Value = 0.2666666666666667
$ws.Range("L12").Value = 0.06666666666666667
# Row 13
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
# Row 15
$ws.Range("F15").Value = 0.02777777777777778
$ws.Range("H15").Value = 0.1388888888888889
$ws.Range("I15").Value = 0.05555555555555555
$ws.Range("J15").Value = 0.3611111111111111
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("O15").Value = 0.1388888888888889
$ws.Range("S15").Value = 0.2222222222222222
# Row 16
$ws.Range("F16").Value = 0.03571428571428571
$ws.Range("H16").Value = 0.2857142857142857
$ws.Range("J16").Value = 0.4642857142857143
$ws.Range("K16").Value = 0.03571428571428571
$ws.Range("O16").Value = 0.1071428571428571
$ws.Range("S16").Value = 0.07142857142857142
# Row 17
$ws.Range("F17").Value = 0.01030927835051546
$ws.Range("H17").Value = 0.1237113402061856
$ws.Range("I17").Value = 0.04123711340206185
$ws.Range("J17").Value = 0.6082474226804123
$ws.Range("K17").Value = 0.09278350515463918
$ws.Range("M17").Value = 0.02061855670103093
$ws.Range("O17").Value = 0.07216494845360824
$ws.Range("S17").Value = 0.03092783505154639
# Row 18
$ws.Range("F18").Value = 0.02222222222222222
$ws.Range("H18").Value = 0.1777777777777778
$ws.Range("I18").Value = 0.08888888888888889
$ws.Range("J18").Value = 0.5555555555555556
$ws.Range("K18").Value = 0.04444444444444445
$ws.Range("O18").Value = 0.02222222222222222
$ws.Range("S18").Value = 0.08888888888888889
# Row 19
$ws.Range("H19").Value = 0.1551724137931035
$ws.Range("I19").Value = 0.09770114942528736
$ws.Range("J19").Value = 0.5574712643678161
$ws.Range("K19").Value = 0.04597701149425287
$ws.Range("M19").Value = 0.005747126436781609
$ws.Range("O19").Value = 0.05172413793103448
$ws.Range("S19").Value = 0.08620689655172414
